$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H40").Value = 4291.3335
$ws_ALC.Range("J40").Value = 4799.8
$ws_ALC.Range("L40").Value = 4799.8
$ws_ALC.Range("N40").Value = -5149.8
$ws_ALC.Range("H44").Value = 5050
$ws_ALC.Range("J44").Value = 5050
$ws_ALC.Range("L44").Value = 5050
$ws_ALC.Range("N44").Value = -5974
$ws_ALC.Range("H69").Value = 12331.667
$ws_ALC.Range("I69").Value = 14995
$ws_ALC.Range("J69").Value = 11000
$ws_ALC.Range("K69").Value = 44985
$ws_ALC.Range("L69").Value = 33000
$ws_ALC.Range("M69").Value = -44111
$ws_ALC.Range("N69").Value = -34748
$ws_ALC.Range("H70").Value = 7981.9287
$ws_ALC.Range("I70").Value = 20000
$ws_ALC.Range("J70").Value = 7057.4614
$ws_ALC.Range("K70").Value = 60000
$ws_ALC.Range("L70").Value = 21172.3842
$ws_ALC.Range("M70").Value = -59730
$ws_ALC.Range("N70").Value = -21712.3842
$ws_ALC.Range("H72").Value = 12331.667
$ws_ALC.Range("I72").Value = 14995
$ws_ALC.Range("J72").Value = 11000
$ws_ALC.Range("K72").Value = 134955
$ws_ALC.Range("L72").Value = 99000
$ws_ALC.Range("M72").Value = -130587
$ws_ALC.Range("N72").Value = -107736
$ws_ALC.Range("H73").Value = 7981.9287
$ws_ALC.Range("I73").Value = 20000
$ws_ALC.Range("J73").Value = 7057.4614
$ws_ALC.Range("K73").Value = 60000
$ws_ALC.Range("L73").Value = 21172.3842
$ws_ALC.Range("M73").Value = -59064
$ws_ALC.Range("N73").Value = -23044.3842
$ws_ALC.Range("H98").Value = 37988.59
$ws_ALC.Range("I98").Value = 55918.4
$ws_ALC.Range("K98").Value = 55918.4
$ws_ALC.Range("M98").Value = -54420.4
$ws_ALC.Range("H107").Value = 10678.5
$ws_ALC.Range("I107").Value = 11643.889
$ws_ALC.Range("K107").Value = 11643.889
$ws_ALC.Range("M107").Value = -9723.888999999999
$ws_ALC.Range("H112").Value = 38309.91
$ws_ALC.Range("J112").Value = 39413.344
$ws_ALC.Range("L112").Value = 118240.032
$ws_ALC.Range("N112").Value = -120456.032
$ws_ALC.Range("H115").Value = 1203
$ws_ALC.Range("I115").Value = 460
$ws_ALC.Range("J115").Value = 1698.3334
$ws_ALC.Range("K115").Value = 1380
$ws_ALC.Range("L115").Value = 5095.0002
$ws_ALC.Range("M115").Value = 187
$ws_ALC.Range("N115").Value = -8229.0002
$ws_ALC.Range("H122").Value = 37988.59
$ws_ALC.Range("I122").Value = 55918.4
$ws_ALC.Range("K122").Value = 167755.2
$ws_ALC.Range("M122").Value = -165305.2
$ws_ALC.Range("H137").Value = 523440.9
$ws_ALC.Range("I137").Value = 1147374.1
$ws_ALC.Range("K137").Value = 3442122.3
$ws_ALC.Range("M137").Value = -3439572.3
$ws_ALC.Range("H138").Value = 6718.0835
$ws_ALC.Range("I138").Value = 2000
$ws_ALC.Range("J138").Value = 7147
$ws_ALC.Range("K138").Value = 6000
$ws_ALC.Range("L138").Value = 21441
$ws_ALC.Range("M138").Value = -860
$ws_ALC.Range("N138").Value = -31721
$ws_ARM.Range("H5").Value = 512.375
$ws_ARM.Range("I5").Value = 512.375
$ws_ARM.Range("J5").Value = 0
$ws_ARM.Range("K5").Value = 512.375
$ws_ARM.Range("L5").Value = 0
$ws_ARM.Range("M5").Value = -400.375
$ws_ARM.Range("N5").ClearContents()
$ws_ARM.Range("H32").Value = 2393.3582
$ws_ARM.Range("I32").Value = 2538.3547
$ws_ARM.Range("K32").Value = 2538.3547
$ws_ARM.Range("M32").Value = -2251.3547
$ws_ARM.Range("H45").Value = 157001.64
$ws_ARM.Range("I45").Value = 217467.4
$ws_ARM.Range("K45").Value = 217467.4
$ws_ARM.Range("M45").Value = -217090.4
$ws_ARM.Range("H88").Value = 1710.8572
$ws_ARM.Range("I88").Value = 1144.375
$ws_ARM.Range("J88").Value = 2466.1667
$ws_ARM.Range("K88").Value = 1144.375
$ws_ARM.Range("L88").Value = 2466.1667
$ws_ARM.Range("M88").Value = -738.375
$ws_ARM.Range("N88").Value = -3278.1667
$ws_ARM.Range("H91").Value = 1710.8572
$ws_ARM.Range("I91").Value = 1144.375
$ws_ARM.Range("J91").Value = 2466.1667
$ws_ARM.Range("K91").Value = 1144.375
$ws_ARM.Range("L91").Value = 2466.1667
$ws_ARM.Range("M91").Value = 259.625
$ws_ARM.Range("N91").Value = -5274.1667
$ws_ARM.Range("H102").Value = 5984.4614
$ws_ARM.Range("I102").Value = 4198.1
$ws_ARM.Range("J102").Value = 7864.8423
$ws_ARM.Range("K102").Value = 4198.1
$ws_ARM.Range("L102").Value = 7864.8423
$ws_ARM.Range("M102").Value = -2576.1
$ws_ARM.Range("N102").Value = -11108.8423
$ws_ARM.Range("H103").Value = 81249.75
$ws_ARM.Range("J103").Value = 81249.75
$ws_ARM.Range("L103").Value = 81249.75
$ws_ARM.Range("N103").Value = -83593.75
$ws_ARM.Range("H106").Value = 0
$ws_ARM.Range("J106").Value = 0
$ws_ARM.Range("L106").Value = 0
$ws_ARM.Range("N106").ClearContents()
$ws_ARM.Range("H111").Value = 89500
$ws_ARM.Range("J111").Value = 89500
$ws_ARM.Range("L111").Value = 89500
$ws_ARM.Range("N111").Value = -97680
$ws_ARM.Range("H122").Value = 641084.9399999999
$ws_ARM.Range("I122").Value = 4201.6924
$ws_ARM.Range("K122").Value = 12605.0772
$ws_ARM.Range("M122").Value = -10155.0772
$ws_BSM.Range("H4").Value = 512.375
$ws_BSM.Range("I4").Value = 512.375
$ws_BSM.Range("J4").Value = 0
$ws_BSM.Range("K4").Value = 512.375
$ws_BSM.Range("L4").Value = 0
$ws_BSM.Range("M4").Value = -397.375
$ws_BSM.Range("N4").ClearContents()
$ws_BSM.Range("H58").Value = 0
$ws_BSM.Range("J58").Value = 0
$ws_BSM.Range("L58").Value = 0
$ws_BSM.Range("N58").ClearContents()
$ws_BSM.Range("H99").Value = 16304.654
$ws_BSM.Range("I99").Value = 19124.904
$ws_BSM.Range("J99").Value = 4459.6
$ws_BSM.Range("K99").Value = 19124.904
$ws_BSM.Range("L99").Value = 4459.6
$ws_BSM.Range("M99").Value = -17626.904
$ws_BSM.Range("N99").Value = -7455.6
$ws_BSM.Range("H123").Value = 0
$ws_BSM.Range("I123").Value = 0
$ws_BSM.Range("J123").Value = 0
$ws_BSM.Range("K123").Value = 0
$ws_BSM.Range("L123").Value = 0
$ws_BSM.Range("M123").ClearContents()
$ws_BSM.Range("N123").ClearContents()
$ws_BSM.Range("H134").Value = 8133.28
$ws_BSM.Range("I134").Value = 8405.781999999999
$ws_BSM.Range("K134").Value = 25217.346
$ws_BSM.Range("M134").Value = -22682.346
$ws_CRP.Range("H31").Value = 3296.5757
$ws_CRP.Range("I31").Value = 1768.7142
$ws_CRP.Range("J31").Value = 3707.923
$ws_CRP.Range("K31").Value = 1768.7142
$ws_CRP.Range("L31").Value = 3707.923
$ws_CRP.Range("M31").Value = -1473.7142
$ws_CRP.Range("N31").Value = -4297.923
$ws_CRP.Range("H34").Value = 3296.5757
$ws_CRP.Range("I34").Value = 1768.7142
$ws_CRP.Range("J34").Value = 3707.923
$ws_CRP.Range("K34").Value = 1768.7142
$ws_CRP.Range("L34").Value = 3707.923
$ws_CRP.Range("M34").Value = -1566.7142
$ws_CRP.Range("N34").Value = -4111.923
$ws_CRP.Range("H43").Value = 19719.6
$ws_CRP.Range("J43").Value = 19719.6
$ws_CRP.Range("L43").Value = 19719.6
$ws_CRP.Range("N43").Value = -20087.6
$ws_CRP.Range("H62").Value = 9146.125
$ws_CRP.Range("I62").Value = 7624.5
$ws_CRP.Range("K62").Value = 7624.5
$ws_CRP.Range("M62").Value = -7000.5
$ws_CRP.Range("H65").Value = 9146.125
$ws_CRP.Range("I65").Value = 7624.5
$ws_CRP.Range("K65").Value = 38122.5
$ws_CRP.Range("M65").Value = -35002.5
$ws_CRP.Range("H101").Value = 19719.6
$ws_CRP.Range("J101").Value = 19719.6
$ws_CRP.Range("L101").Value = 19719.6
$ws_CRP.Range("N101").Value = -26209.6
$ws_CRP.Range("H134").Value = 1567749.8
$ws_CRP.Range("I134").Value = 2610272
$ws_CRP.Range("K134").Value = 7830816
$ws_CRP.Range("M134").Value = -7828281
$ws_CUL.Range("H68").Value = 6040.6113
$ws_CUL.Range("J68").Value = 7355.2593
$ws_CUL.Range("L68").Value = 22065.7779
$ws_CUL.Range("N68").Value = -23687.7779
$ws_CUL.Range("H71").Value = 6040.6113
$ws_CUL.Range("J71").Value = 7355.2593
$ws_CUL.Range("L71").Value = 66197.3337
$ws_CUL.Range("N71").Value = -74309.3337
$ws_CUL.Range("H131").Value = 5850.375
$ws_CUL.Range("J131").Value = 2762.75
$ws_CUL.Range("L131").Value = 8288.25
$ws_CUL.Range("N131").Value = -18368.25
$ws_CUL.Range("H133").Value = 9909.6
$ws_CUL.Range("I133").Value = 9909.6
$ws_CUL.Range("J133").Value = 0
$ws_CUL.Range("K133").Value = 29728.8
$ws_CUL.Range("L133").Value = 0
$ws_CUL.Range("M133").Value = -24668.8
$ws_CUL.Range("N133").ClearContents()
$ws_LTW.Range("H46").Value = 3110.9546
$ws_LTW.Range("I46").Value = 1220
$ws_LTW.Range("J46").Value = 3667.1177
$ws_LTW.Range("K46").Value = 1220
$ws_LTW.Range("L46").Value = 3667.1177
$ws_LTW.Range("M46").Value = -1032
$ws_LTW.Range("N46").Value = -4043.1177
$ws_LTW.Range("H100").Value = 8045.5835
$ws_LTW.Range("I100").Value = 8904.700000000001
$ws_LTW.Range("J100").Value = 3750
$ws_LTW.Range("K100").Value = 8904.700000000001
$ws_LTW.Range("L100").Value = 3750
$ws_LTW.Range("M100").Value = -8363.700000000001
$ws_LTW.Range("N100").Value = -4832
$ws_LTW.Range("H132").Value = 23686.385
$ws_LTW.Range("I132").Value = 31770.334
$ws_LTW.Range("K132").Value = 95311.00199999999
$ws_LTW.Range("M132").Value = -92781.00199999999
$ws_LTW.Range("H136").Value = 5426.1333
$ws_LTW.Range("I136").Value = 4873.3335
$ws_LTW.Range("J136").Value = 5564.3335
$ws_LTW.Range("K136").Value = 14620.0005
$ws_LTW.Range("L136").Value = 16693.0005
$ws_LTW.Range("M136").Value = -12070.0005
$ws_LTW.Range("N136").Value = -21793.0005
$ws_WVR.Range("H62").Value = 103207.734
$ws_WVR.Range("I62").Value = 148784.5
$ws_WVR.Range("K62").Value = 148784.5
$ws_WVR.Range("M62").Value = -148160.5
$ws_WVR.Range("H65").Value = 103207.734
$ws_WVR.Range("I65").Value = 148784.5
$ws_WVR.Range("K65").Value = 743922.5
$ws_WVR.Range("M65").Value = -740802.5
$ws_WVR.Range("H81").Value = 30320.715
$ws_WVR.Range("I81").Value = 34586.668
$ws_WVR.Range("J81").Value = 4725
$ws_WVR.Range("K81").Value = 69173.336
$ws_WVR.Range("L81").Value = 9450
$ws_WVR.Range("M81").Value = -68112.336
$ws_WVR.Range("N81").Value = -11572
$ws_WVR.Range("H84").Value = 30320.715
$ws_WVR.Range("I84").Value = 34586.668
$ws_WVR.Range("J84").Value = 4725
$ws_WVR.Range("K84").Value = 345866.68
$ws_WVR.Range("L84").Value = 47250
$ws_WVR.Range("M84").Value = -340562.68
$ws_WVR.Range("N84").Value = -57858
$ws_WVR.Range("H132").Value = 30437.6
$ws_WVR.Range("I132").Value = 34961.273
$ws_WVR.Range("J132").Value = 17997.5
$ws_WVR.Range("K132").Value = 104883.819
$ws_WVR.Range("L132").Value = 53992.5
$ws_WVR.Range("M132").Value = -102353.819
$ws_WVR.Range("N132").Value = -59052.5
$ws_WVR.Range("H136").Value = 4371.143
$ws_WVR.Range("I136").Value = 3433
$ws_WVR.Range("K136").Value = 10299
$ws_WVR.Range("M136").Value = -7749
